$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 3.3
$ws.Range("L2").Value = 3.75
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.62
$ws.Range("R2").Value = 2.3
$ws.Range("S2").Value = 2.05
$ws.Range("T2").Value = 1.85
$ws.Range("W2").Value = 1.3
$ws.Range("X2").Value = 3.4
$ws.Range("Y2").Value = 1.53
$ws.Range("Z2").Value = 2.38
$ws.Range("AA2").Value = 11
$ws.Range("AC2").Value = 9.5
$ws.Range("AF2").Value = 23
$ws.Range("AK2").Value = 151

# Row 3
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 2

# Row 4
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 2.63
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 2.05
$ws.Range("R4").Value = 1.85
$ws.Range("U4").Value = 3.5
$ws.Range("V4").Value = 1.3
$ws.Range("AA4").Value = 9.5
$ws.Range("AK4").Value = 301
$ws.Range("AL4").Value = 9.5

# Row 5
$ws.Range("N5").Value = 29

# Row 6
$ws.Range("I6").Value = 3.5
$ws.Range("K6").Value = 2.15
$ws.Range("L6").Value = 3.95
$ws.Range("N6").Value = 7.8
$ws.Range("W6").Value = 1.38
$ws.Range("X6").Value = 2.82
$ws.Range("Z6").Value = 2.12
$ws.Range("AA6").Value = 8.5
$ws.Range("AC6").Value = 8.25
$ws.Range("AG6").Value = 7.8
$ws.Range("AL6").Value = 11.25
$ws.Range("AP6").Value = 30
$ws.Range("AQ6").Value = 35

# Row 7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
$ws.Range("W7").Value = 1.54
$ws.Range("AR7").Value = 1.95
$ws.Range("AS7").Value = 1.9

# Row 8
$ws.Range("W8").Value = 1.41
$ws.Range("X8").Value = 2.62

# Row 9
$ws.Range("G9").Value = 2.9
$ws.Range("I9").Value = 2.35
$ws.Range("J9").Value = 3.5
$ws.Range("W9").Value = 1.37
$ws.Range("AD9").Value = 29
$ws.Range("AM9").Value = 12
$ws.Range("AO9").Value = 23

# Row 10
$ws.Range("W10").Value = 1.41
$ws.Range("X10").Value = 2.62

# Row 11
$ws.Range("W11").Value = 1.37

# Row 12
$ws.Range("G12").Value = 2.9
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("W12").Value = 1.54
$ws.Range("AC12").Value = 12
$ws.Range("AD12").Value = 34
$ws.Range("AP12").Value = 23

# Row 13
$ws.Range("G13").Value = 1.7
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 2.3
$ws.Range("K13").Value = 2.4
$ws.Range("Q13").Value = 1.6
$ws.Range("R13").Value = 2.3
$ws.Range("S13").Value = 1.98
$ws.Range("T13").Value = 1.83
$ws.Range("W13").Value = 1.27
$ws.Range("Y13").Value = 1.57
$ws.Range("Z13").Value = 2.25
$ws.Range("AH13").Value = 8

# Row 14
$ws.Range("G14").Value = 1.39
$ws.Range("J14").Value = 1.87
$ws.Range("M14").Value = 1.02
$ws.Range("O14").Value = 1.15
$ws.Range("V14").Value = 1.47

# Row 15
$ws.Range("G15").Value = 4.5
$ws.Range("H15").Value = 4.1
$ws.Range("I15").Value = 1.58
$ws.Range("J15").Value = 4.75
$ws.Range("L15").Value = 2.1
$ws.Range("M15").Value = 1.02
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.17
$ws.Range("U15").Value = 2.62
$ws.Range("V15").Value = 1.41
$ws.Range("W15").Value = 1.33
$ws.Range("X15").Value = 3.25
$ws.Range("Y15").Value = 1.73
$ws.Range("Z15").Value = 2
$ws.Range("AF15").Value = 41
$ws.Range("AM15").Value = 8.5
$ws.Range("AO15").Value = 12

# Row 16
$ws.Range("G16").Value = 1.1
$ws.Range("J16").Value = 1.37
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("O16").Value = 1.05
$ws.Range("P16").Value = 7.5
$ws.Range("Q16").Value = 1.25
$ws.Range("R16").Value = 3.75
$ws.Range("U16").Value = 1.63
$ws.Range("V16").Value = 2.1

# Row 17
$ws.Range("G17").Value = 5.5
$ws.Range("I17").Value = 1.37
$ws.Range("K17").Value = 2.6
$ws.Range("L17").Value = 1.8
$ws.Range("M17").Value = 21
$ws.Range("N17").Value = 1.03
$ws.Range("O17").Value = 1.11
$ws.Range("Q17").Value = 1.44
$ws.Range("R17").Value = 2.63
$ws.Range("U17").Value = 2.1
$ws.Range("V17").Value = 1.63
$ws.Range("AA17").Value = 21
$ws.Range("AC17").Value = 19
$ws.Range("AD17").Value = 67
$ws.Range("AG17").Value = 21
$ws.Range("AN17").Value = 9
$ws.Range("AO17").Value = 10

# Row 18
$ws.Range("M18").Value = 1.03
$ws.Range("O18").Value = 1.22
$ws.Range("V18").Value = 1.33

# Row 19
$ws.Range("G19").Value = 1.11
$ws.Range("O19").Value = 1.05
$ws.Range("U19").Value = 1.76
$ws.Range("V19").Value = 1.96

# Row 20
$ws.Range("M20").Value = 1.03
$ws.Range("O20").Value = 1.25
$ws.Range("V20").Value = 1.3

# Row 21
$ws.Range("M21").Value = 1.05
$ws.Range("O21").Value = 1.3
$ws.Range("Q21").Value = 2.08
$ws.Range("R21").Value = 1.73
$ws.Range("V21").Value = 1.22

# Row 22
$ws.Range("G22").Value = 1.41
$ws.Range("M22").Value = 1.03
$ws.Range("O22").Value = 1.27
$ws.Range("V22").Value = 1.27

# Row 23
$ws.Range("G23").Value = 3.8
$ws.Range("H23").Value = 3.55
$ws.Range("I23").Value = 1.9
$ws.Range("J23").Value = 4.15
$ws.Range("K23").Value = 2.15
$ws.Range("L23").Value = 2.47
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 8
$ws.Range("O23").Value = 1.27
$ws.Range("P23").Value = 3.5
$ws.Range("Q23").Value = 1.82
$ws.Range("R23").Value = 1.93
$ws.Range("U23").Value = 2.95
$ws.Range("V23").Value = 1.37
$ws.Range("W23").Value = 1.4
$ws.Range("X23").Value = 2.82
$ws.Range("Z23").Value = 2.02
$ws.Range("AA23").Value = 11.25
$ws.Range("AC23").Value = 13.5
$ws.Range("AF23").Value = 40
$ws.Range("AG23").Value = 8
$ws.Range("AH23").Value = 7.2
$ws.Range("AI23").Value = 15
$ws.Range("AK23").Value = 500
$ws.Range("AL23").Value = 7.5
$ws.Range("AQ23").Value = 27
